$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet ("NC1 (Tiện ĐK trong)" -> "NGOẠI QUAN")
$ws.Name = "NGOẠI QUAN"

# 2. Update the Print_Area defined name so it points at the new sheet name
foreach ($n in $wb.Names) {
    if ($n.Name -like "*Print_Area*") {
        $n.RefersTo = "='NGOẠI QUAN'!`$B`$1:`$N`$22"
    }
}

# 3. Update the report title in B1
$ws.Range("B1").Value = "THEO DÕI TIẾN ĐỘ KẾ HOẠCH NGOẠI QUAN VBAT5A-9-12"

# 4. Update the start/end plan dates
$ws.Range("D4").Value = 45749.51783564815
$ws.Range("D5").Value = 45758.532476851855

# 5. Update the first (still relevant) tracking row - row 9
$ws.Range("B9").Value = 45758.51783564815
$ws.Range("C9").Value = "VBAT5A-9-12"
$ws.Range("D9").Value = "NGOẠI QUAN"
$ws.Range("E9").Value = "B1"
$ws.Range("G9").Value = ""
$ws.Range("H9").Value = ""
$ws.Range("I9").Value = "Hoàng Thuận"
$ws.Range("J9").Value = 397
$ws.Range("K9").Value = 397
$ws.Range("M9").Value = 0

# 6. Remove (clear) the no-longer-used tracking rows 10-17 while keeping
#    their existing cell formatting intact.
$ws.Range("B10:N17").ClearContents()
